# Update column G ("K") values on Sheet1 to reflect regenerated save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 6
    3  = 3
    4  = 3
    5  = 0
    6  = 6
    7  = 7
    8  = 4
    9  = 3
    10 = 6
    11 = 4
    12 = 10
    13 = 6
    14 = 2
    15 = 10
    16 = 2
    17 = 5
    18 = 2
    19 = 5
    20 = 2
    21 = 2
    22 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
